$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 681
$ws1.Range("F4").Value = 935
$ws1.Range("F5").Value = 705
$ws1.Range("F6").Value = 831
$ws1.Range("F7").Value = 392
$ws1.Range("F8").Value = 592
$ws1.Range("F9").Value = 125
$ws1.Range("F10").Value = 1194
$ws1.Range("F12").Value = 375
$ws1.Range("F13").Value = 497
$ws1.Range("F16").Value = 382
$ws1.Range("F17").Value = 339
$ws1.Range("F20").Value = 551
$ws1.Range("F21").Value = 65
$ws1.Range("F22").Value = 565
$ws1.Range("F23").Value = 25
$ws1.Range("F24").Value = 706

# ---- Sheet: 演出 (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G3").Value = 138
$ws2.Range("F4").Value = 312
$ws2.Range("F8").Value = 176
$ws2.Range("F11").Value = 23
$ws2.Range("F13").Value = 71

# ---- Sheet: 本地生活 (Local life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 367

# ---- Sheet: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 367
$ws4.Range("G5").Value = 138
$ws4.Range("F6").Value = 312
$ws4.Range("F7").Value = 681
$ws4.Range("F8").Value = 935
$ws4.Range("F9").Value = 705
$ws4.Range("F10").Value = 831
$ws4.Range("F11").Value = 392
$ws4.Range("F12").Value = 592
$ws4.Range("F13").Value = 125
$ws4.Range("F14").Value = 1194
$ws4.Range("F18").Value = 375
$ws4.Range("F19").Value = 497
$ws4.Range("F23").Value = 382
$ws4.Range("F24").Value = 176
$ws4.Range("F25").Value = 339
$ws4.Range("F30").Value = 551
$ws4.Range("F31").Value = 23
$ws4.Range("F33").Value = 71
$ws4.Range("F34").Value = 65
$ws4.Range("F35").Value = 565
$ws4.Range("F36").Value = 25
$ws4.Range("F37").Value = 706
